# Project Update Till 100424
# - Remove the 4 trailing (now-unused) rows from "Problem Solving Competency"
#   and clear the leftover stray 0/placeholder values that preceded them.
# - Clear the stale placeholder numbers in "Monthly Data" (B2:C13).
# - Reset each sheet's view: drop the old ad-hoc selections/scroll positions
#   and freeze row 1 with the cursor parked at A2, making
#   "Problem Solving Competency" the active sheet/tab again.

$wb = $excel.ActiveWorkbook

# --- Sheet: Problem Solving Competency ---
$ws1 = $wb.Worksheets.Item("Problem Solving Competency")

# Drop the last 4 rows (345:348) - data no longer tracked, shift remaining rows up.
$ws1.Rows("345:348").Delete() | Out-Null

# Rows 319:344 kept their dates but the metric columns are cleared out (no data yet).
$ws1.Range("B319:F344").ClearContents()

# --- Sheet: Weekly Data ---
$ws2 = $wb.Worksheets.Item("Weekly Data")
$ws2.Activate()
$ws2.Rows.Item(2).Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- Sheet: Monthly Data ---
$ws3 = $wb.Worksheets.Item("Monthly Data")
$ws3.Range("B2:C13").ClearContents()
$ws3.Activate()
$ws3.Rows.Item(2).Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- Re-activate Problem Solving Competency sheet & set its view/selection ---
$ws1.Activate()
$ws1.Rows.Item(2).Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
